$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C24").Value = "ок"
$ws.Range("D24").Value = "ок"
$ws.Range("E24").Value = "ок"

Write-Host "done"
